# Auto-generated Excel COM-interop script applying the Masamune_Profits diff
# Updates computed profit/price columns (H-N) across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 4680.5
$ws.Range("I62").Value = 5973.636
$ws.Range("J62").Value = 3100
$ws.Range("K62").Value = 5973.636
$ws.Range("L62").Value = 3100
$ws.Range("M62").Value = -5349.636
$ws.Range("N62").Value = -4348
$ws.Range("H65").Value = 4680.5
$ws.Range("I65").Value = 5973.636
$ws.Range("J65").Value = 3100
$ws.Range("K65").Value = 29868.18
$ws.Range("L65").Value = 15500
$ws.Range("M65").Value = -26748.18
$ws.Range("N65").Value = -21740
$ws.Range("H98").Value = 26126.156
$ws.Range("I98").Value = 1552.5
$ws.Range("K98").Value = 1552.5
$ws.Range("M98").Value = -54.5
$ws.Range("H99").Value = 960.05554
$ws.Range("I99").Value = 1093.25
$ws.Range("J99").Value = 693.6667
$ws.Range("K99").Value = 3279.75
$ws.Range("L99").Value = 2081.0001
$ws.Range("M99").Value = -1781.75
$ws.Range("N99").Value = -5077.0001
$ws.Range("H101").Value = 824.875
$ws.Range("I101").Value = 603.1667
$ws.Range("J101").Value = 1490
$ws.Range("K101").Value = 1809.5001
$ws.Range("L101").Value = 4470
$ws.Range("M101").Value = -187.5001
$ws.Range("N101").Value = -7714
$ws.Range("H122").Value = 26126.156
$ws.Range("I122").Value = 1552.5
$ws.Range("K122").Value = 4657.5
$ws.Range("M122").Value = -2207.5
$ws.Range("H129").Value = 1054.95
$ws.Range("J129").Value = 968.3469
$ws.Range("L129").Value = 2905.0407
$ws.Range("N129").Value = -12905.0407
$ws.Range("H137").Value = 1511322.1
$ws.Range("I137").Value = 2331837
$ws.Range("K137").Value = 6995511
$ws.Range("M137").Value = -6992961

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1705.3462
$ws.Range("I2").Value = 1771.875
$ws.Range("J2").Value = 907
$ws.Range("K2").Value = 1771.875
$ws.Range("L2").Value = 907
$ws.Range("M2").Value = -1658.875
$ws.Range("N2").Value = -1133
$ws.Range("H32").Value = 11375.143
$ws.Range("I32").Value = 10530.898
$ws.Range("K32").Value = 10530.898
$ws.Range("M32").Value = -10243.898
$ws.Range("H82").Value = 10000
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 10000
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
$ws.Range("H116").Value = 1705.3462
$ws.Range("I116").Value = 1771.875
$ws.Range("J116").Value = 907
$ws.Range("K116").Value = 1771.875
$ws.Range("L116").Value = 907
$ws.Range("M116").Value = 522.125
$ws.Range("N116").Value = -5495
$ws.Range("H122").Value = 1752.2632
$ws.Range("I122").Value = 1727.5333
$ws.Range("J122").Value = 1845
$ws.Range("K122").Value = 5182.5999
$ws.Range("L122").Value = 5535
$ws.Range("M122").Value = -2732.5999
$ws.Range("N122").Value = -10435

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1705.3462
$ws.Range("I3").Value = 1771.875
$ws.Range("J3").Value = 907
$ws.Range("K3").Value = 1771.875
$ws.Range("L3").Value = 907
$ws.Range("M3").Value = -1657.875
$ws.Range("N3").Value = -1135
$ws.Range("H10").Value = 70006
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("H94").Value = 813.4545000000001
$ws.Range("I94").Value = 852.2857
$ws.Range("K94").Value = 852.2857
$ws.Range("M94").Value = -401.2857
$ws.Range("H134").Value = 2636.6753
$ws.Range("I134").Value = 1552.8611
$ws.Range("J134").Value = 3588.3171
$ws.Range("K134").Value = 4658.5833
$ws.Range("L134").Value = 10764.9513
$ws.Range("M134").Value = -2123.5833
$ws.Range("N134").Value = -15834.9513

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2366.16
$ws.Range("I31").Value = 823.5263
$ws.Range("J31").Value = 3311.6453
$ws.Range("K31").Value = 823.5263
$ws.Range("L31").Value = 3311.6453
$ws.Range("M31").Value = -528.5263
$ws.Range("N31").Value = -3901.6453
$ws.Range("H33").Value = 3200
$ws.Range("I33").Value = 2900
$ws.Range("K33").Value = 2900
$ws.Range("M33").Value = -2521
$ws.Range("H34").Value = 2366.16
$ws.Range("I34").Value = 823.5263
$ws.Range("J34").Value = 3311.6453
$ws.Range("K34").Value = 823.5263
$ws.Range("L34").Value = 3311.6453
$ws.Range("M34").Value = -621.5263
$ws.Range("N34").Value = -3715.6453

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3521.0557
$ws.Range("I5").Value = 4466.72
$ws.Range("J5").Value = 1371.8182
$ws.Range("K5").Value = 13400.16
$ws.Range("L5").Value = 4115.4546
$ws.Range("M5").Value = -13288.16
$ws.Range("N5").Value = -4339.4546
$ws.Range("H12").Value = 326.96667
$ws.Range("I12").Value = 248.28572
$ws.Range("J12").Value = 350.91306
$ws.Range("K12").Value = 744.85716
$ws.Range("L12").Value = 1052.73918
$ws.Range("M12").Value = -571.85716
$ws.Range("N12").Value = -1398.73918
$ws.Range("H68").Value = 3426.4546
$ws.Range("I68").Value = 694
$ws.Range("J68").Value = 4033.6667
$ws.Range("K68").Value = 2082
$ws.Range("L68").Value = 12101.0001
$ws.Range("M68").Value = -1271
$ws.Range("N68").Value = -13723.0001
$ws.Range("H71").Value = 3426.4546
$ws.Range("I71").Value = 694
$ws.Range("J71").Value = 4033.6667
$ws.Range("K71").Value = 6246
$ws.Range("L71").Value = 36303.0003
$ws.Range("M71").Value = -2190
$ws.Range("N71").Value = -44415.0003
$ws.Range("H107").Value = 11477.167
$ws.Range("I107").Value = 12921.125
$ws.Range("J107").Value = 10322
$ws.Range("K107").Value = 38763.375
$ws.Range("L107").Value = 30966
$ws.Range("M107").Value = -36843.375
$ws.Range("N107").Value = -34806
$ws.Range("H113").Value = 5243.409
$ws.Range("I113").Value = 8770.083000000001
$ws.Range("J113").Value = 1011.4
$ws.Range("K113").Value = 26310.249
$ws.Range("L113").Value = 3034.2
$ws.Range("M113").Value = -24140.249
$ws.Range("N113").Value = -7374.2
$ws.Range("H131").Value = 971.6799999999999
$ws.Range("I131").Value = 466.66666
$ws.Range("J131").Value = 987.29895
$ws.Range("K131").Value = 1399.99998
$ws.Range("L131").Value = 2961.89685
$ws.Range("M131").Value = 3640.00002
$ws.Range("N131").Value = -13041.89685
$ws.Range("H135").Value = 3521.0557
$ws.Range("I135").Value = 4466.72
$ws.Range("J135").Value = 1371.8182
$ws.Range("K135").Value = 40200.48
$ws.Range("L135").Value = 12346.3638
$ws.Range("M135").Value = -37665.48
$ws.Range("N135").Value = -17416.3638

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 825
$ws.Range("I36").Value = 800
$ws.Range("J36").Value = 850
$ws.Range("K36").Value = 800
$ws.Range("L36").Value = 850
$ws.Range("M36").Value = -315
$ws.Range("N36").Value = -1820
$ws.Range("H43").Value = 299701.88
$ws.Range("I43").Value = 1668311.4
$ws.Range("J43").Value = 6428.4287
$ws.Range("K43").Value = 1668311.4
$ws.Range("L43").Value = 6428.4287
$ws.Range("M43").Value = -1668160.4
$ws.Range("N43").Value = -6730.4287
$ws.Range("H46").Value = 27426.334
$ws.Range("J46").Value = 27426.334
$ws.Range("L46").Value = 27426.334
$ws.Range("N46").Value = -27738.334
$ws.Range("H112").Value = 20000
$ws.Range("J112").Value = 20000
$ws.Range("L112").Value = 20000
$ws.Range("N112").Value = -22216
$ws.Range("H122").Value = 1531.4166
$ws.Range("I122").Value = 1195.4
$ws.Range("J122").Value = 1771.4286
$ws.Range("K122").Value = 3586.2
$ws.Range("L122").Value = 5314.2858
$ws.Range("M122").Value = -1136.2
$ws.Range("N122").Value = -10214.2858

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H111").Value = 46000
$ws.Range("J111").Value = 46000
$ws.Range("L111").Value = 46000
$ws.Range("N111").Value = -54180
$ws.Range("H140").Value = 23499.5
$ws.Range("J140").Value = 23499.5
$ws.Range("L140").Value = 23499.5
$ws.Range("N140").Value = -33859.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1361434.6
$ws.Range("I122").Value = 2198649.2
$ws.Range("J122").Value = 960.625
$ws.Range("K122").Value = 6595947.600000001
$ws.Range("L122").Value = 2881.875
$ws.Range("M122").Value = -6593497.600000001
$ws.Range("N122").Value = -7781.875
$ws.Range("H126").Value = 1401764.9
$ws.Range("I126").Value = 1731344.9
$ws.Range("J126").Value = 1049.75
$ws.Range("K126").Value = 5194034.699999999
$ws.Range("L126").Value = 3149.25
$ws.Range("M126").Value = -5191564.699999999
$ws.Range("N126").Value = -8089.25
$ws.Range("H132").Value = 1318850.1
$ws.Range("I132").Value = 1500291.6
$ws.Range("K132").Value = 4500874.800000001
$ws.Range("M132").Value = -4500874.800000001
